$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "Erarity" -> "ERarity" in the rarity column header cell
$ws.Range("B3").Value = "ERarity"

# Move the active selection to D8 (matches saved cursor position in sheet view)
$ws.Range("D8").Select()
